# Add 4 new "Whitebox" test rows (T016-T019) for the validDestination
# function, just below the existing "getInput" test table (before the
# "checkBoxSize" section). This shifts every following section down by
# 3 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 already exists (blank spacer row); insert 3 more rows after it
# so the new test data occupies rows 24-27 and the next section header
# (previously row 25) lands on row 28.
$ws.Rows("25:27").Insert()

# Fill in the new cells. The assignment order below matches the order
# in which these distinct strings were first typed into the workbook,
# so that shared-string allocation mirrors the original edit.
$ws.Range("A24").Value = 'T016'
$ws.Range("B27").Value = 'WhiteBox'
$ws.Range("C27").Value = 'valid destination in input string'
$ws.Range("E27").Value = '0 0 1A'
$ws.Range("F24").Value = 'pass'
$ws.Range("F25").Value = 'pass'
$ws.Range("F26").Value = 'pass'
$ws.Range("F27").Value = 'pass'
$ws.Range("G27").Value = 'Function correctly detects correct destination input'
$ws.Range("A25").Value = 'T017'
$ws.Range("A26").Value = 'T018'
$ws.Range("B24").Value = 'Whitebox'
$ws.Range("B25").Value = 'Whitebox'
$ws.Range("B26").Value = 'Whitebox'
$ws.Range("C26").Value = 'no map, null map provided'
$ws.Range("D26").Value = 'struct Map map* = NULL'
$ws.Range("G26").Value = 'Function correctly detects there is no map provided'
$ws.Range("A27").Value = 'T019'
$ws.Range("C24").Value = 'no input, null input provided'
$ws.Range("D24").Value = 'char* input = NULL'
$ws.Range("E24").Value = '"\0"'
$ws.Range("G24").Value = 'Function correctly detects null or empty input string'
$ws.Range("C25").Value = 'exit String'
$ws.Range("D25").Value = 'char input[] = "0 0 x"'
$ws.Range("D27").Value = 'char input[] = "0 0 1A"'
$ws.Range("E26").Value = 'NULL'
$ws.Range("G25").Value = 'Function corrects detects input string and returns corresponding return value'
$ws.Range("E25").Value = '0 0 x'

# Match the page orientation change recorded for this sheet.
$ws.PageSetup.Orientation = 1

# Restore the selection/active cell noted for this edit.
$ws.Range("G25").Select()
